$wb = $excel.ActiveWorkbook

# --- Sheet "Lab3Rubric_CS295N" (sheet1) ---
$ws1 = $wb.Worksheets.Item("Lab3Rubric_CS295N")

# Update the point values for the "Data base stores:" section
$ws1.Range("B11:C11").Value = 10
$ws1.Range("B12:C12").Value = 10
$ws1.Range("B13:C13").Value = 8
$ws1.Range("B14:C14").Value = 7

# Apply the new font color/style to B11:C14 (creates the new font + cell style)
$ws1.Range("B11:C14").Font.Color = 0

# Update the subtotal formulas to reference the corrected range
$ws1.Range("B15").Formula = "=SUM(B11:B14)"
$ws1.Range("C15").Formula = "=SUM(C11:C14)"

# Update sheet selection
[void]$ws1.Range("F22").Select()

# --- Sheet "Student Points" (sheet2) ---
$ws2 = $wb.Worksheets.Item("Student Points")

$ws2.Range("B11:C11").Value = 10
$ws2.Range("B12:C12").Value = 10
$ws2.Range("B13:C13").Value = 8
$ws2.Range("B14:C14").Value = 7

$ws2.Range("B15").Formula = "=SUM(B11:B14)"
$ws2.Range("C15").Formula = "=SUM(C11:C14)"

[void]$ws2.Range("E12").Select()

# Make "Student Points" the active (visible) tab
[void]$ws2.Activate()
